# "Fix capital letter bug during test generation"
#
# The measurement-collection script used a case-sensitive match (e.g.
# "casc" vs "CASC"/mis-cased config flag) that made it write the
# "direct-XSTS, casc, ..." rows into the wrong block of the sheet: the
# values that belong in the A98:A107 / A110:A119 block ended up in the
# A225:A234 / A237:A246 block instead (and vice versa). The fix simply
# swaps the two blocks of raw measurements back to where they belong;
# the dependent AVERAGE/MEDIAN formulas recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Block 1: A98:A107  <->  A225:A234  (10 raw samples each)
$block1Top    = $ws.Range("A98:A107").Value2
$block1Bottom = $ws.Range("A225:A234").Value2
$ws.Range("A98:A107").Value2   = $block1Bottom
$ws.Range("A225:A234").Value2 = $block1Top

# Block 2: A110:A119  <->  A237:A246  (10 raw samples each)
$block2Top    = $ws.Range("A110:A119").Value2
$block2Bottom = $ws.Range("A237:A246").Value2
$ws.Range("A110:A119").Value2   = $block2Bottom
$ws.Range("A237:A246").Value2 = $block2Top

# Leave the sheet scrolled/selected where the fix was made.
$ws.Range("A97:E119").Select()
